# Assignment6 Client Command Table - update error-checking notes to mention
# saving the HTM service/characteristic *handle* (not just the service/char
# itself), per "change error checking to LOG_ERROR" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two notes cells with the refined wording.
$ws.Range("E24").Value = " save HTM service handle to data structure for rreference "
$ws.Range("E26").Value = "save to HTM characteristic handle to structure for reference "

# Reflect the author's updated view state: zoomed to 90%, with F26 selected.
$excel.ActiveWindow.Zoom = 90
$ws.Range("F26").Select()
